$wb = $excel.ActiveWorkbook

# --- Sheet "life table" (1st sheet): fix merge-cell bookkeeping order ---
# (no content change on this sheet other than the merged-range registration
# order, which Excel rewrites when the ranges are re-merged)
$wsLifeTable = $wb.Worksheets.Item(1)
$wsLifeTable.Range("I12:I16").UnMerge()
$wsLifeTable.Range("I12:I16").Merge()
$wsLifeTable.Range("A17:A21").UnMerge()
$wsLifeTable.Range("A17:A21").Merge()
$wsLifeTable.Range("A22:A26").UnMerge()
$wsLifeTable.Range("A22:A26").Merge()
$wsLifeTable.Range("A27:A32").UnMerge()
$wsLifeTable.Range("A27:A32").Merge()
$wsLifeTable.Range("I27:I32").UnMerge()
$wsLifeTable.Range("I27:I32").Merge()
$wsLifeTable.Range("I17:I21").UnMerge()
$wsLifeTable.Range("I17:I21").Merge()
$wsLifeTable.Range("I22:I26").UnMerge()
$wsLifeTable.Range("I22:I26").Merge()

# --- Sheet "Dispersal" (2nd sheet): fix add_agent bug so first-gen males
#     are added to the male set (probability values updated) ---
$wsDispersal = $wb.Worksheets.Item(2)

# Update probabilities: column F (0.5 -> 0.75) for rows 7-32
$wsDispersal.Range("F7:F32").Value = 0.75

# Update probabilities: column E (0.4 -> 0.5) for rows 9-32
$wsDispersal.Range("E9:E32").Value = 0.5

# Re-register merged ranges so they end up ordered after the ones below
# (matches the order Excel produces when the sheet is resaved)
$wsDispersal.Range("L34:M34").UnMerge()
$wsDispersal.Range("L34:M34").Merge()
$wsDispersal.Range("A22:A26").UnMerge()
$wsDispersal.Range("A22:A26").Merge()
$wsDispersal.Range("A27:A32").UnMerge()
$wsDispersal.Range("A27:A32").Merge()
$wsDispersal.Range("A4:A6").UnMerge()
$wsDispersal.Range("A4:A6").Merge()
$wsDispersal.Range("A7:A8").UnMerge()
$wsDispersal.Range("A7:A8").Merge()
$wsDispersal.Range("A9:A11").UnMerge()
$wsDispersal.Range("A9:A11").Merge()
$wsDispersal.Range("B33:H33").UnMerge()
$wsDispersal.Range("B33:H33").Merge()

# Update the sheet view (scroll position + selection)
$wsDispersal.Activate()
$wsDispersal.Range("A10").Select()
$excel.ActiveWindow.ScrollRow = 10
$wsDispersal.Range("E6:E32").Select()
